$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-29 Sunday" "2024-12-30 Monday"

Replace-Text "627÷4=" "392÷3="
Replace-Text "732÷9=" "692÷2="
Replace-Text "115÷3=" "313÷3="
Replace-Text "646÷7=" "606÷8="
Replace-Text "283÷8=" "504÷8="

Replace-Text "326÷6=" "201÷2="
Replace-Text "124÷4=" "502÷7="
Replace-Text "416÷6=" "815÷4="
Replace-Text "318÷7=" "721÷3="
Replace-Text "890÷9=" "946÷6="

Replace-Text "775÷5=" "598÷4="
Replace-Text "249÷7=" "908÷2="
Replace-Text "919÷3=" "871÷6="
Replace-Text "886÷5=" "674÷6="
Replace-Text "728÷7=" "320÷9="

Replace-Text "379÷9=" "575÷6="
Replace-Text "880÷3=" "563÷3="
Replace-Text "982÷3=" "644÷4="
Replace-Text "518÷7=" "482÷8="
Replace-Text "303÷6=" "262÷3="

Replace-Text "614÷8=" "331÷6="
Replace-Text "822÷9=" "739÷5="
Replace-Text "154÷8=" "492÷2="
Replace-Text "597÷6=" "316÷2="
Replace-Text "551÷5=" "463÷8="
